$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table (and fix the swapped
# Filecoin/ONDO rows) per the Wed Jun 19 18:33:39 UTC 2024 data pull.
# Column D values are numeric-looking strings (e.g. "65.076.56"); a
# leading apostrophe forces Excel to keep them as text instead of
# auto-converting them to numbers, matching the original text cells.
$ws.Range("D2").Value = "'65.076.56"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "'3.564.58"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'600.15"
$ws.Range("E5").Value = "  +3.19%  "
$ws.Range("D6").Value = "'137.96"
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("D7").Value = "'3.563.27"
$ws.Range("E7").Value = "  +4.42%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +3.35%  "
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("E12").Value = "  +4.22%  "
$ws.Range("D13").Value = "'4.180.37"
$ws.Range("E13").Value = "  +4.69%  "
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("E15").Value = "  +5.16%  "
$ws.Range("D16").Value = "'3.579.37"
$ws.Range("E16").Value = "  +4.93%  "
$ws.Range("D17").Value = "'0.116"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "'64.999.59"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "'10.13"
$ws.Range("E19").Value = "  +8.35%  "
$ws.Range("D20").Value = "'14.42"
$ws.Range("E20").Value = "  +7.32%  "
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("D22").Value = "'391.84"
$ws.Range("E22").Value = "  +3.11%  "
$ws.Range("D23").Value = "'0.579"
$ws.Range("E23").Value = "  +7.30%  "
$ws.Range("D24").Value = "'3.710.65"
$ws.Range("E24").Value = "  +4.48%  "
$ws.Range("D25").Value = "'74.17"
$ws.Range("E25").Value = "  +3.51%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'0.0000117"
$ws.Range("E27").Value = "  +13.64%  "
$ws.Range("D28").Value = "'7.65"
$ws.Range("E28").Value = "  +6.35%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +5.35%  "
$ws.Range("D31").Value = "'8.34"
$ws.Range("E31").Value = "  +4.83%  "
$ws.Range("D32").Value = "'3.576.83"
$ws.Range("E32").Value = "  +4.18%  "
$ws.Range("E33").Value = "  +21.32%  "
$ws.Range("D34").Value = "'24.04"
$ws.Range("E34").Value = "  +5.33%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").Value = "'170.64"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  +5.22%  "
$ws.Range("D39").Value = "'1.55"
$ws.Range("E39").Value = "  +7.87%  "
$ws.Range("D40").Value = "'5.01"
$ws.Range("E40").Value = "  +9.94%  "
$ws.Range("D41").Value = "'0.0809"
$ws.Range("E41").Value = "  +6.97%  "
$ws.Range("D42").Value = "'0.830"
$ws.Range("E42").Value = "  +3.83%  "
$ws.Range("D43").Value = "'26.88"
$ws.Range("E43").Value = "  +21.24%  "
$ws.Range("D44").Value = "'42.66"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'4.47"
$ws.Range("E46").Value = "  +5.21%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.22"
$ws.Range("E47").Value = "  +9.81%  "
$ws.Range("D48").Value = "'1.66"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("D49").Value = "'2.466.11"
$ws.Range("E49").Value = "  +12.43%  "
$ws.Range("E50").Value = "  +7.04%  "
$ws.Range("D51").Value = "'2.36"
$ws.Range("E51").Value = "  +16.24%  "
